$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.745631333333333
$ws.Range("H2").Value = 11.236894
$ws.Range("I2").Value = 0.1419671142338921
$ws.Range("J2").Value = 0.1419671142338921
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.06089466666666667
$ws.Range("N2").Value = 0.182684
$ws.Range("O2").Value = 0.001903591634475228
$ws.Range("P2").Value = 0.001903591634475228
$ws.Range("Q2").Value = 0.2280889714995556
$ws.Range("R2").Value = 2.052800743496
$ws.Range("S2").Value = 0.000270247411026226
$ws.Range("T2").Value = 0.000270247411026226
# Row 3
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.745631333333333
$ws.Range("H3").Value = 11.236894
$ws.Range("I3").Value = 0.1419671142338921
$ws.Range("J3").Value = 0.1419671142338921
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 28.046323
$ws.Range("N3").Value = 84.138969
$ws.Range("O3").Value = 0.8767392739472014
$ws.Range("P3").Value = 0.8767392739472013
$ws.Range("Q3").Value = 105.0511862135873
$ws.Range("R3").Value = 945.460675922286
$ws.Range("S3").Value = 0.1244681446578019
$ws.Range("T3").Value = 0.1244681446578019
# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.745631333333333
$ws.Range("H4").Value = 11.236894
$ws.Range("I4").Value = 0.1419671142338921
$ws.Range("J4").Value = 0.1419671142338921
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 3.882136333333333
$ws.Range("N4").Value = 11.646409
$ws.Range("O4").Value = 0.1213571344183235
$ws.Range("P4").Value = 0.1213571344183235
$ws.Range("Q4").Value = 14.54105149040511
$ws.Range("R4").Value = 130.869463413646
$ws.Range("S4").Value = 0.01722872216506392
$ws.Range("T4").Value = 0.01722872216506392
# Row 5
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 15.68955666666667
$ws.Range("H5").Value = 47.06867
$ws.Range("I5").Value = 0.5946663954227359
$ws.Range("J5").Value = 0.5946663954227359
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.06089466666666667
$ws.Range("N5").Value = 0.182684
$ws.Range("O5").Value = 0.001903591634475228
$ws.Range("P5").Value = 0.001903591634475228
$ws.Range("Q5").Value = 0.9554103233644445
$ws.Range("R5").Value = 8.59869291028
$ws.Range("S5").Value = 0.001132001975630258
$ws.Range("T5").Value = 0.001132001975630258
# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 15.68955666666667
$ws.Range("H6").Value = 47.06867
$ws.Range("I6").Value = 0.5946663954227359
$ws.Range("J6").Value = 0.5946663954227359
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 28.046323
$ws.Range("N6").Value = 84.138969
$ws.Range("O6").Value = 0.8767392739472014
$ws.Range("P6").Value = 0.8767392739472013
$ws.Range("Q6").Value = 440.0343740001367
$ws.Range("R6").Value = 3960.30936600123
$ws.Range("S6").Value = 0.5213673837637288
$ws.Range("T6").Value = 0.5213673837637288
# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 15.68955666666667
$ws.Range("H7").Value = 47.06867
$ws.Range("I7").Value = 0.5946663954227359
$ws.Range("J7").Value = 0.5946663954227359
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 3.882136333333333
$ws.Range("N7").Value = 11.646409
$ws.Range("O7").Value = 0.1213571344183235
$ws.Range("P7").Value = 0.1213571344183235
$ws.Range("Q7").Value = 60.90899798955889
$ws.Range("R7").Value = 548.18098190603
$ws.Range("S7").Value = 0.07216700968337686
$ws.Range("T7").Value = 0.07216700968337686
# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Tfpi"
$ws.Range("C8").Value = "Vldlr"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 6.948607666666667
$ws.Range("H8").Value = 20.845823
$ws.Range("I8").Value = 0.263366490343372
$ws.Range("J8").Value = 0.263366490343372
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.06089466666666667
$ws.Range("N8").Value = 0.182684
$ws.Range("O8").Value = 0.001903591634475228
$ws.Range("P8").Value = 0.001903591634475228
$ws.Range("Q8").Value = 0.4231331476591111
$ws.Range("R8").Value = 3.808198328932
$ws.Range("S8").Value = 0.0005013422478187439
$ws.Range("T8").Value = 0.0005013422478187438
# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Tfpi"
$ws.Range("C9").Value = "Vldlr"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 6.948607666666667
$ws.Range("H9").Value = 20.845823
$ws.Range("I9").Value = 0.263366490343372
$ws.Range("J9").Value = 0.263366490343372
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 28.046323
$ws.Range("N9").Value = 84.138969
$ws.Range("O9").Value = 0.8767392739472014
$ws.Range("P9").Value = 0.8767392739472013
$ws.Range("Q9").Value = 194.8828950196097
$ws.Range("R9").Value = 1753.946055176487
$ws.Range("S9").Value = 0.2309037455256706
$ws.Range("T9").Value = 0.2309037455256706
# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Tfpi"
$ws.Range("C10").Value = "Vldlr"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 6.948607666666667
$ws.Range("H10").Value = 20.845823
$ws.Range("I10").Value = 0.263366490343372
$ws.Range("J10").Value = 0.263366490343372
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.882136333333333
$ws.Range("N10").Value = 11.646409
$ws.Range("O10").Value = 0.1213571344183235
$ws.Range("P10").Value = 0.1213571344183235
$ws.Range("Q10").Value = 26.97544228884522
$ws.Range("R10").Value = 242.778980599607
$ws.Range("S10").Value = 0.03196140256988268
$ws.Range("T10").Value = 0.03196140256988268
